# Weekly update: shift historical Albahaca price records for
# "Vega Modelo de Temuco" down by one row (rows 373-412 each take on
# the prior content of the row above them) and seed row 372 with the
# newest week's record (new date + quantity, same price/unit/region).
# This mirrors the upstream weekly commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 372
$ws.Range("D372").Value = 45180
$ws.Range("J372").Value = 90
# Row 373
$ws.Range("D373").Value = 44810
$ws.Range("J373").Value = 80
$ws.Range("K373").Value = 6000
$ws.Range("L373").Value = 6000
$ws.Range("M373").Value = 6000
$ws.Range("P373").Value = 6000
# Row 374
$ws.Range("D374").Value = 44340
# Row 375
$ws.Range("J375").Value = 50
$ws.Range("O375").Value = "Región de Arica y Parinacota"
# Row 376
$ws.Range("D376").Value = 44175
$ws.Range("J376").Value = 80
$ws.Range("K376").Value = 5000
$ws.Range("L376").Value = 5000
$ws.Range("M376").Value = 5000
$ws.Range("O376").Value = "Región del Maule"
$ws.Range("P376").Value = 5000
# Row 377
$ws.Range("D377").Value = 44746
$ws.Range("J377").Value = 55
$ws.Range("K377").Value = 6000
$ws.Range("L377").Value = 6000
$ws.Range("M377").Value = 6000
$ws.Range("O377").Value = "Región de Arica y Parinacota"
$ws.Range("P377").Value = 6000
# Row 378
$ws.Range("J378").Value = 30
$ws.Range("O378").Value = "Región de La Araucanía"
# Row 379
$ws.Range("D379").Value = 44273
$ws.Range("J379").Value = 80
$ws.Range("K379").Value = 5000
$ws.Range("L379").Value = 5000
$ws.Range("M379").Value = 5000
$ws.Range("O379").Value = "Región del Maule"
$ws.Range("P379").Value = 5000
# Row 380
$ws.Range("D380").Value = 44799
$ws.Range("J380").Value = 50
# Row 381
$ws.Range("D381").Value = 45093
$ws.Range("J381").Value = 45
# Row 382
$ws.Range("D382").Value = 44426
$ws.Range("J382").Value = 40
$ws.Range("K382").Value = 6000
$ws.Range("L382").Value = 6000
$ws.Range("M382").Value = 6000
$ws.Range("P382").Value = 6000
# Row 383
$ws.Range("D383").Value = 44181
$ws.Range("J383").Value = 35
$ws.Range("K383").Value = 5000
$ws.Range("L383").Value = 5000
$ws.Range("M383").Value = 5000
$ws.Range("P383").Value = 5000
# Row 384
$ws.Range("D384").Value = 45075
$ws.Range("J384").Value = 65
$ws.Range("O384").Value = "Región de Arica y Parinacota"
# Row 385
$ws.Range("J385").Value = 50
$ws.Range("K385").Value = 6000
$ws.Range("L385").Value = 6000
$ws.Range("M385").Value = 6000
$ws.Range("O385").Value = "Región de La Araucanía"
$ws.Range("P385").Value = 6000
# Row 386
$ws.Range("D386").Value = 44970
$ws.Range("J386").Value = 100
$ws.Range("K386").Value = 4000
$ws.Range("L386").Value = 4000
$ws.Range("M386").Value = 4000
$ws.Range("O386").Value = "Región del Maule"
$ws.Range("P386").Value = 4000
# Row 387
$ws.Range("J387").Value = 30
$ws.Range("K387").Value = 5000
$ws.Range("L387").Value = 5000
$ws.Range("M387").Value = 5000
$ws.Range("O387").Value = "Región de La Araucanía"
$ws.Range("P387").Value = 5000
# Row 388
$ws.Range("D388").Value = 44960
$ws.Range("J388").Value = 100
$ws.Range("K388").Value = 4000
$ws.Range("L388").Value = 4000
$ws.Range("M388").Value = 4000
$ws.Range("O388").Value = "Región del Maule"
$ws.Range("P388").Value = 4000
# Row 389
$ws.Range("D389").Value = 44883
$ws.Range("J389").Value = 50
$ws.Range("K389").Value = 8000
$ws.Range("L389").Value = 8000
$ws.Range("M389").Value = 8000
$ws.Range("O389").Value = "Región Metropolitana"
$ws.Range("P389").Value = 8000
# Row 390
$ws.Range("J390").Value = 125
$ws.Range("K390").Value = 6000
$ws.Range("L390").Value = 6000
$ws.Range("M390").Value = 6000
$ws.Range("O390").Value = "Región de La Araucanía"
$ws.Range("P390").Value = 6000
# Row 391
$ws.Range("D391").Value = 44588
$ws.Range("J391").Value = 65
$ws.Range("K391").Value = 5000
$ws.Range("L391").Value = 5000
$ws.Range("M391").Value = 5000
$ws.Range("O391").Value = "Región del Maule"
$ws.Range("P391").Value = 5000
# Row 392
$ws.Range("D392").Value = 44749
$ws.Range("J392").Value = 80
$ws.Range("K392").Value = 6000
$ws.Range("L392").Value = 6000
$ws.Range("M392").Value = 6000
$ws.Range("O392").Value = "Región de Arica y Parinacota"
$ws.Range("P392").Value = 6000
# Row 393
$ws.Range("J393").Value = 30
$ws.Range("K393").Value = 7000
$ws.Range("L393").Value = 7000
$ws.Range("M393").Value = 7000
$ws.Range("O393").Value = "Región de La Araucanía"
$ws.Range("P393").Value = 7000
# Row 394
$ws.Range("D394").Value = 44579
$ws.Range("J394").Value = 50
$ws.Range("K394").Value = 5000
$ws.Range("L394").Value = 5000
$ws.Range("M394").Value = 5000
$ws.Range("O394").Value = "Región del Maule"
$ws.Range("P394").Value = 5000
# Row 395
$ws.Range("D395").Value = 44413
$ws.Range("J395").Value = 30
$ws.Range("K395").Value = 8000
$ws.Range("L395").Value = 8000
$ws.Range("M395").Value = 8000
$ws.Range("P395").Value = 8000
# Row 396
$ws.Range("D396").Value = 44351
$ws.Range("J396").Value = 15
$ws.Range("K396").Value = 5000
$ws.Range("L396").Value = 5000
$ws.Range("M396").Value = 5000
$ws.Range("O396").Value = "Región de Arica y Parinacota"
$ws.Range("P396").Value = 5000
# Row 397
$ws.Range("D397").Value = 44544
$ws.Range("J397").Value = 75
$ws.Range("K397").Value = 6000
$ws.Range("L397").Value = 7000
$ws.Range("M397").Value = 6467
$ws.Range("O397").Value = "Región del Maule"
$ws.Range("P397").Value = 6467
# Row 398
$ws.Range("D398").Value = 44453
$ws.Range("J398").Value = 20
$ws.Range("O398").Value = "Región de Arica y Parinacota"
# Row 399
$ws.Range("D399").Value = 44901
$ws.Range("J399").Value = 35
$ws.Range("K399").Value = 8000
$ws.Range("L399").Value = 8000
$ws.Range("M399").Value = 8000
$ws.Range("O399").Value = "Región Metropolitana"
$ws.Range("P399").Value = 8000
# Row 400
$ws.Range("J400").Value = 50
$ws.Range("O400").Value = "Región de La Araucanía"
# Row 401
$ws.Range("D401").Value = 44217
$ws.Range("J401").Value = 80
$ws.Range("O401").Value = "Región del Maule"
# Row 402
$ws.Range("D402").Value = 44767
$ws.Range("J402").Value = 100
$ws.Range("K402").Value = 5000
$ws.Range("L402").Value = 5000
$ws.Range("M402").Value = 5000
$ws.Range("O402").Value = "Región de Arica y Parinacota"
$ws.Range("P402").Value = 5000
# Row 403
$ws.Range("D403").Value = 44599
$ws.Range("J403").Value = 65
$ws.Range("K403").Value = 7000
$ws.Range("L403").Value = 7000
$ws.Range("M403").Value = 7000
$ws.Range("O403").Value = "Región del Maule"
$ws.Range("P403").Value = 7000
# Row 404
$ws.Range("D404").Value = 44880
$ws.Range("J404").Value = 30
$ws.Range("K404").Value = 8000
$ws.Range("L404").Value = 8000
$ws.Range("M404").Value = 8000
$ws.Range("O404").Value = "Región Metropolitana"
$ws.Range("P404").Value = 8000
# Row 405
$ws.Range("D405").Value = 45117
$ws.Range("J405").Value = 90
$ws.Range("K405").Value = 5000
$ws.Range("L405").Value = 5000
$ws.Range("M405").Value = 5000
$ws.Range("O405").Value = "Región de Arica y Parinacota"
$ws.Range("P405").Value = 5000
# Row 406
$ws.Range("D406").Value = 44637
$ws.Range("J406").Value = 80
$ws.Range("K406").Value = 7000
$ws.Range("L406").Value = 7000
$ws.Range("M406").Value = 7000
$ws.Range("O406").Value = "Región de La Araucanía"
$ws.Range("P406").Value = 7000
# Row 407
$ws.Range("D407").Value = 44362
$ws.Range("J407").Value = 25
$ws.Range("K407").Value = 5000
$ws.Range("L407").Value = 5000
$ws.Range("M407").Value = 5000
$ws.Range("O407").Value = "Región de Arica y Parinacota"
$ws.Range("P407").Value = 5000
# Row 408
$ws.Range("J408").Value = 50
$ws.Range("O408").Value = "Región Metropolitana"
# Row 409
$ws.Range("D409").Value = 44893
$ws.Range("J409").Value = 40
$ws.Range("K409").Value = 9000
$ws.Range("L409").Value = 9000
$ws.Range("M409").Value = 9000
$ws.Range("P409").Value = 9000
# Row 410
$ws.Range("D410").Value = 44557
$ws.Range("J410").Value = 55
$ws.Range("L410").Value = 5000
$ws.Range("M410").Value = 5000
$ws.Range("O410").Value = "Región del Maule"
$ws.Range("P410").Value = 5000
# Row 411
$ws.Range("J411").Value = 90
$ws.Range("L411").Value = 6000
$ws.Range("M411").Value = 5389
$ws.Range("O411").Value = "Región de La Araucanía"
$ws.Range("P411").Value = 5389
# Row 412
$ws.Range("D412").Value = 44242
$ws.Range("K412").Value = 5000
$ws.Range("L412").Value = 5000
$ws.Range("M412").Value = 5000
$ws.Range("N412").Value = "`$/paquete"
$ws.Range("O412").Value = "Región del Maule"
$ws.Range("P412").Value = 5000
$ws.Range("Q412").Value = 1
